$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (becomes the old row 7 data - "Blåsippa")
$ws.Range("A4").Value = 111634290
$ws.Range("B4").Value = 98535
$ws.Range("E4").Value = 222498
$ws.Range("F4").Value = "Blåsippa"
$ws.Range("G4").Value = "Hepatica nobilis"
$ws.Range("H4").Value = "Schreb."

# Row 5 (becomes the old row 4 data - "Fjällig taggsvamp s.str.")
$ws.Range("A5").Value = 111634304
$ws.Range("B5").Value = 90687
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 5964
$ws.Range("F5").Value = "Fjällig taggsvamp s.str."
$ws.Range("G5").Value = "Sarcodon imbricatus s.str."
$ws.Range("H5").Value = "(L.:Fr.) P.Karst."
$ws.Range("Q5").Value = 676708.8668162767
$ws.Range("R5").Value = 6618511.450801066

# Row 7 (becomes the old row 5 data - "Orange taggsvamp")
$ws.Range("A7").Value = 111633890
$ws.Range("B7").Value = 90658
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 4361
$ws.Range("F7").Value = "Orange taggsvamp"
$ws.Range("G7").Value = "Hydnellum aurantiacum"
$ws.Range("H7").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q7").Value = 676486.710397501
$ws.Range("R7").Value = 6618439.724061669
